# "Titulo de validar por Excel"
# Adds a second column (B) to the "articulo" sheet holding the shopping-cart
# page title: B1 = "tituloPgCarro" (header), B2 = "Carrito de compras" (value).
# Also widens column B and keeps column A's text formatting consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("articulo")
$ws.Activate()

# New header + value in column B
$ws.Range("B1").Value = "tituloPgCarro"
$ws.Range("B2").Value = "Carrito de compras"

# Keep/ensure text ("@") number format across the now 2-column data block
$ws.Range("A1:B2").NumberFormat = "@"

# Give the new column a sensible custom width (matches ~17.36 chars on the
# sheet's pixel grid)
$ws.Columns.Item(2).ColumnWidth = 16.5

# Move the active selection, as recorded after the edit
$ws.Range("C8").Select() | Out-Null
